$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Calculate the percentage of cash and online payment"
#    -> "Count number of cash and online payments:"
#    (also re-homes the _GoBack bookmark to just after this text)
# ------------------------------------------------------------------

$found = $d.Content.Find.Execute("Calculate the percentage of cash and online payment", `
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

$target = $d.Content
$target.Find.Execute("Calculate the percentage of cash and online payment", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null

$r = $target
$r.Text = "Count number of cash and online payments:X"

# Move the (hidden) _GoBack bookmark to sit right before the "X" marker,
# which is currently a non-boundary (mid-paragraph) position.
$bmPos = $r.End - 1
$rBm = $d.Range($bmPos, $bmPos)

$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()
$d.Bookmarks.Add("_GoBack", $rBm)

# Now drop the temporary "X" marker - this leaves the bookmark exactly at
# the end of the paragraph's text (immediately before the paragraph mark).
$rX = $d.Range($r.End - 1, $r.End)
$rX.Delete()

# ------------------------------------------------------------------
# 2) Merge '... = "' / 'online' / '" Reservation)' into one run of text
#    (a plain text replace collapses them back into a single run).
# ------------------------------------------------------------------

$d.Content.Find.Execute([char]8220 + "online" + [char]8221 + " Reservation)", $true, $false, $false, $false, $false, $true, 1, $false, `
    [char]8220 + "online" + [char]8221 + " Reservation)", 2) | Out-Null
